$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-16 Monday" "2024-09-17 Tuesday"

Replace-Text "299×7=" "195×7="
Replace-Text "531×4=" "552×8="
Replace-Text "765×7=" "709×9="
Replace-Text "759×9=" "248×6="
Replace-Text "531×6=" "234×7="
Replace-Text "459×6=" "305×9="
Replace-Text "683×3=" "260×4="
Replace-Text "687×4=" "409×3="
Replace-Text "826×5=" "314×7="
Replace-Text "978×2=" "149×4="
Replace-Text "827×3=" "243×5="
Replace-Text "933×7=" "386×2="
Replace-Text "402×7=" "182×7="
Replace-Text "357×7=" "976×5="
Replace-Text "285×4=" "778×5="
Replace-Text "884×4=" "656×2="
Replace-Text "693×3=" "554×3="
Replace-Text "946×8=" "120×4="
Replace-Text "344×5=" "589×7="
Replace-Text "560×4=" "930×4="
Replace-Text "601×8=" "856×2="
Replace-Text "901×3=" "457×4="
Replace-Text "494×7=" "514×8="
Replace-Text "921×4=" "613×6="
Replace-Text "150×5=" "187×4="
